# Auto-generated PowerShell Excel COM-interop script
# Applies: insert 'fertilized' & 'hatchery_fry_release' columns/rows
$wb = $excel.ActiveWorkbook

# ---------- Sheet1 (metadata): insert 2 rows at row 12 ----------
$ws1 = $wb.Worksheets.Item("metadata")
$ws1.Range("A12:A13").EntireRow.Insert()

$ws1.Range("A12").Value = "fertilized"
$ws1.Range("B12").Value = "Binary variable describing whether (1) or not (0) the CU nursary lake was fertilized in each year. Note that fertilization affects abundances of pre-smolts in year + 1 (e.g. fertilizing a lake in 2010 is expected to bolster the abundances of pre-smolts counted during the winter 2011 ATS."
$ws1.Range("A13").Value = "hatchery_fry_release"
$ws1.Range("B13").Value = "Numbers of Sockeye fry released by the Hucuktlis/Henderson hatchery in each year. Pertains only to the Hucuktlis CU because outplants were discontinued in the 1930s. The hatchery ceased operations in 2007."

# ---------- Sheet2 (S-R data): insert 2 columns at K:L ----------
$ws2 = $wb.Worksheets.Item("S-R data")
$ws2.Range("K1:L1").EntireColumn.Insert()

$ws2.Range("K1").Value = "fertilized"
$ws2.Range("L1").Value = "hatchery_fry_release"

# ---------- Fill in 'fertilized' (K) values per row ----------
$fertilized = @{
    2 = 1
    3 = 0
    4 = 1
    5 = 0
    6 = 1
    7 = 0
    8 = 1
    9 = 0
    10 = 1
    11 = 0
    12 = 1
    13 = 0
    14 = 1
    15 = 0
    16 = 1
    17 = 0
    18 = 1
    19 = 1
    20 = 1
    21 = 0
    22 = 1
    23 = 0
    24 = 1
    25 = 0
    26 = 1
    27 = 0
    28 = 1
    29 = 0
    30 = 1
    31 = 0
    32 = 1
    33 = 0
    34 = 1
    35 = 0
    36 = 1
    37 = 0
    38 = 1
    39 = 0
    40 = 1
    41 = 0
    42 = 1
    43 = 0
    44 = 1
    45 = 0
    46 = 1
    47 = 0
    48 = 1
    49 = 0
    50 = 1
    51 = 0
    52 = 1
    53 = 0
    54 = 1
    55 = 0
    56 = 1
    57 = 0
    58 = 1
    59 = 0
    60 = 1
    61 = 0
    62 = 1
    63 = 0
    64 = 1
    65 = 0
    66 = 1
    67 = 0
    68 = 1
    69 = 0
    70 = 1
    71 = 0
    72 = 1
    73 = 0
    74 = 1
    75 = 0
    76 = 1
    77 = 0
    78 = 1
    79 = 0
    80 = 1
    81 = 0
    82 = 1
    83 = 0
    84 = 1
    85 = 0
    86 = 1
    87 = 0
    88 = 1
    89 = 0
    90 = 1
    91 = 0
    92 = 1
    93 = 0
    94 = 1
    95 = 0
    96 = 0
    97 = 0
    98 = 0
    99 = 0
    100 = 1
    101 = 1
    102 = 1
    103 = 1
    104 = 1
    105 = 1
    106 = 1
    107 = 1
    108 = 1
    109 = 1
    110 = 1
    111 = 1
    112 = 1
    113 = 1
    114 = 1
    115 = 1
    116 = 1
    117 = 1
    118 = 1
    119 = 1
    120 = 1
    121 = 1
    122 = 0
    123 = 1
    124 = 0
    125 = 0
    126 = 0
    127 = 0
    128 = 1
    129 = 0
    130 = 0
    131 = 0
    132 = 0
    133 = 0
    134 = 0
    135 = 0
    136 = 0
    137 = 0
    138 = 0
    139 = 0
    140 = 0
    141 = 0
    142 = 0
    143 = 0
    144 = 0
}
foreach ($row in $fertilized.Keys) {
    $ws2.Cells.Item([int]$row, 11).Value = $fertilized[$row]
}

# ---------- Fill in 'hatchery_fry_release' (L) values per row ----------
$hatcheryFry = @{
    96 = 0
    97 = 0
    98 = 0
    99 = 0
    100 = 0
    101 = 0
    102 = 0
    103 = 0
    104 = 0
    105 = 0
    106 = 0
    107 = 0
    108 = 0
    109 = 0
    110 = 0
    111 = 0
    112 = 0
    113 = 0
    114 = 0
    115 = 0
    116 = 70000
    117 = 659000
    118 = 658000
    119 = 206000
    120 = 862000
    121 = 1025000
    122 = 860000
    123 = 1200000
    124 = 1900000
    125 = 2100000
    126 = 2300000
    127 = 783000
    128 = 0
    129 = 0
    130 = 0
    131 = 0
    132 = 0
    133 = 0
    134 = 0
    135 = 0
    136 = 0
    137 = 0
    138 = 0
    139 = 0
    140 = 0
    141 = 0
    142 = 0
    143 = 0
    144 = 0
}
foreach ($row in $hatcheryFry.Keys) {
    $ws2.Cells.Item([int]$row, 12).Value = $hatcheryFry[$row]
}

# ---------- Minor floating point correction on F19 ----------
$ws2.Range("F19").Value = 39662.13671849497

Write-Output "done"